$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1899.8334
$ws.Range("I40").Value = 1879.8
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1879.8
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1704.8
$ws.Range("N40").Value = -2350

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H80").Value = 276.66666
$ws.Range("I80").Value = 315.18182
$ws.Range("J80").Value = 216.14285
$ws.Range("K80").Value = 945.54546
$ws.Range("L80").Value = 648.4285500000001
$ws.Range("M80").Value = 52.45453999999995
$ws.Range("N80").Value = -2644.42855

$ws.Range("H83").Value = 276.66666
$ws.Range("I83").Value = 315.18182
$ws.Range("J83").Value = 216.14285
$ws.Range("K83").Value = 2836.63638
$ws.Range("L83").Value = 1945.28565
$ws.Range("M83").Value = 2155.36362
$ws.Range("N83").Value = -11929.28565

$ws.Range("H98").Value = 1500.75
$ws.Range("I98").Value = 1500.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1500.75
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -2.75

$ws.Range("H113").Value = 5719.875
$ws.Range("I113").Value = 5152.4
$ws.Range("J113").Value = 6665.6665
$ws.Range("K113").Value = 5152.4
$ws.Range("L113").Value = 6665.6665
$ws.Range("M113").Value = -1898.4
$ws.Range("N113").Value = -13173.6665

$ws.Range("H122").Value = 1500.75
$ws.Range("I122").Value = 1500.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4502.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2052.25

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H135").Value = 1345.909
$ws.Range("I135").Value = 972.7143
$ws.Range("J135").Value = 1999
$ws.Range("K135").Value = 8754.4287
$ws.Range("L135").Value = 17991
$ws.Range("M135").Value = -6219.4287
$ws.Range("N135").Value = -23061

$ws.Range("H137").Value = 8874
$ws.Range("I137").Value = 2549.9
$ws.Range("J137").Value = 14623.182
$ws.Range("K137").Value = 7649.700000000001
$ws.Range("L137").Value = 43869.546
$ws.Range("M137").Value = -5099.700000000001
$ws.Range("N137").Value = -48969.546

$ws.Range("H138").Value = 5682.716
$ws.Range("I138").Value = 7568.8887
$ws.Range("J138").Value = 5143.8096
$ws.Range("K138").Value = 22706.6661
$ws.Range("L138").Value = 15431.4288
$ws.Range("M138").Value = -17566.6661
$ws.Range("N138").Value = -25711.4288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18185.47
$ws.Range("I32").Value = 8131.3516
$ws.Range("J32").Value = 30185.549
$ws.Range("K32").Value = 8131.3516
$ws.Range("L32").Value = 30185.549
$ws.Range("M32").Value = -7844.3516
$ws.Range("N32").Value = -30759.549

$ws.Range("H45").Value = 2552.875
$ws.Range("I45").Value = 1605.75
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 1605.75
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -1228.75
$ws.Range("N45").Value = -4254

$ws.Range("H61").Value = 2470.6667
$ws.Range("I61").Value = 2320.8572
$ws.Range("J61").Value = 2995
$ws.Range("K61").Value = 2320.8572
$ws.Range("L61").Value = 2995
$ws.Range("M61").Value = -2108.8572
$ws.Range("N61").Value = -3419

$ws.Range("H95").Value = 52999.332
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 52999.332
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 52999.332
$ws.Range("N95").Value = -58491.332

$ws.Range("H97").Value = 593.8182
$ws.Range("I97").Value = 585.2778
$ws.Range("J97").Value = 632.25
$ws.Range("K97").Value = 585.2778
$ws.Range("L97").Value = 632.25
$ws.Range("M97").Value = -89.27779999999996
$ws.Range("N97").Value = -1624.25

$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws.Range("H136").Value = 2470.6667
$ws.Range("I136").Value = 2320.8572
$ws.Range("J136").Value = 2995
$ws.Range("K136").Value = 6962.571599999999
$ws.Range("L136").Value = 8985
$ws.Range("M136").Value = -4412.571599999999
$ws.Range("N136").Value = -14085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2986.0952
$ws.Range("I20").Value = 1681.7693
$ws.Range("J20").Value = 5105.625
$ws.Range("K20").Value = 1681.7693
$ws.Range("L20").Value = 5105.625
$ws.Range("M20").Value = -1434.7693
$ws.Range("N20").Value = -5599.625

$ws.Range("H86").Value = 4123.8
$ws.Range("I86").Value = 3652.5
$ws.Range("J86").Value = 4438
$ws.Range("K86").Value = 3652.5
$ws.Range("L86").Value = 4438
$ws.Range("M86").Value = -2529.5
$ws.Range("N86").Value = -6684

$ws.Range("H89").Value = 4123.8
$ws.Range("I89").Value = 3652.5
$ws.Range("J89").Value = 4438
$ws.Range("K89").Value = 18262.5
$ws.Range("L89").Value = 22190
$ws.Range("M89").Value = -12646.5
$ws.Range("N89").Value = -33422

$ws.Range("H94").Value = 1003.3889
$ws.Range("I94").Value = 972.5625
$ws.Range("J94").Value = 1250
$ws.Range("K94").Value = 972.5625
$ws.Range("L94").Value = 1250
$ws.Range("M94").Value = -521.5625
$ws.Range("N94").Value = -2152

$ws.Range("H105").Value = 3068.468
$ws.Range("I105").Value = 2504.7188
$ws.Range("J105").Value = 4271.1333
$ws.Range("K105").Value = 2504.7188
$ws.Range("L105").Value = 4271.1333
$ws.Range("M105").Value = -757.7188000000001
$ws.Range("N105").Value = -7765.1333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4637.2896
$ws.Range("I31").Value = 3840.3333
$ws.Range("J31").Value = 5621.7646
$ws.Range("K31").Value = 3840.3333
$ws.Range("L31").Value = 5621.7646
$ws.Range("M31").Value = -3545.3333
$ws.Range("N31").Value = -6211.7646

$ws.Range("H34").Value = 4637.2896
$ws.Range("I34").Value = 3840.3333
$ws.Range("J34").Value = 5621.7646
$ws.Range("K34").Value = 3840.3333
$ws.Range("L34").Value = 5621.7646
$ws.Range("M34").Value = -3638.3333
$ws.Range("N34").Value = -6025.7646

$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H62").Value = 73332.836
$ws.Range("I62").Value = 12499.5
$ws.Range("J62").Value = 103749.5
$ws.Range("K62").Value = 12499.5
$ws.Range("L62").Value = 103749.5
$ws.Range("M62").Value = -11875.5
$ws.Range("N62").Value = -104997.5

$ws.Range("H65").Value = 73332.836
$ws.Range("I65").Value = 12499.5
$ws.Range("J65").Value = 103749.5
$ws.Range("K65").Value = 62497.5
$ws.Range("L65").Value = 518747.5
$ws.Range("M65").Value = -59377.5
$ws.Range("N65").Value = -524987.5

$ws.Range("H99").Value = 12374.5
$ws.Range("I99").Value = 9921.182000000001
$ws.Range("J99").Value = 14827.818
$ws.Range("K99").Value = 9921.182000000001
$ws.Range("L99").Value = 14827.818
$ws.Range("M99").Value = -8423.182000000001
$ws.Range("N99").Value = -17823.818

$ws.Range("H105").Value = 4741.3
$ws.Range("I105").Value = 3827.3333
$ws.Range("J105").Value = 5133
$ws.Range("K105").Value = 3827.3333
$ws.Range("L105").Value = 5133
$ws.Range("M105").Value = -2080.3333
$ws.Range("N105").Value = -8627

$ws.Range("H126").Value = 12374.5
$ws.Range("I126").Value = 9921.182000000001
$ws.Range("J126").Value = 14827.818
$ws.Range("K126").Value = 29763.546
$ws.Range("L126").Value = 44483.454
$ws.Range("M126").Value = -27293.546
$ws.Range("N126").Value = -49423.454

$ws.Range("H134").Value = 3296.1052
$ws.Range("I134").Value = 2769.3076
$ws.Range("J134").Value = 4437.5
$ws.Range("K134").Value = 8307.9228
$ws.Range("L134").Value = 13312.5
$ws.Range("M134").Value = -5772.9228
$ws.Range("N134").Value = -18382.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 209.33333
$ws.Range("I23").Value = 214.6
$ws.Range("J23").Value = 207.3077
$ws.Range("K23").Value = 643.8
$ws.Range("L23").Value = 621.9231
$ws.Range("M23").Value = -408.8
$ws.Range("N23").Value = -1091.9231

$ws.Range("H98").Value = 307.75
$ws.Range("I98").Value = 362.4
$ws.Range("J98").Value = 216.66667
$ws.Range("K98").Value = 1087.2
$ws.Range("L98").Value = 650.00001
$ws.Range("M98").Value = 410.8000000000002
$ws.Range("N98").Value = -3646.00001

$ws.Range("H113").Value = 3612.5715
$ws.Range("I113").Value = 2899.3333
$ws.Range("J113").Value = 4147.5
$ws.Range("K113").Value = 8697.999899999999
$ws.Range("L113").Value = 12442.5
$ws.Range("M113").Value = -6527.999899999999
$ws.Range("N113").Value = -16782.5

$ws.Range("H132").Value = 5311.125
$ws.Range("I132").Value = 3250
$ws.Range("J132").Value = 5998.1665
$ws.Range("K132").Value = 29250
$ws.Range("L132").Value = 53983.4985
$ws.Range("M132").Value = -26720
$ws.Range("N132").Value = -59043.4985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H122").Value = 791606.6
$ws.Range("I122").Value = 117722.664
$ws.Range("J122").Value = 2004597.8
$ws.Range("K122").Value = 353167.992
$ws.Range("L122").Value = 6013793.4
$ws.Range("M122").Value = -350717.992
$ws.Range("N122").Value = -6018693.4

$ws.Range("H126").Value = 4996
$ws.Range("I126").Value = 4991
$ws.Range("J126").Value = 4998.5
$ws.Range("K126").Value = 14973
$ws.Range("L126").Value = 14995.5
$ws.Range("M126").Value = -12503
$ws.Range("N126").Value = -19935.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2017.25
$ws.Range("I16").Value = 2070.7
$ws.Range("J16").Value = 1750
$ws.Range("K16").Value = 2070.7
$ws.Range("L16").Value = 1750
$ws.Range("M16").Value = -1900.7
$ws.Range("N16").Value = -2090

$ws.Range("H22").Value = 4449.5
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4449.5
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4449.5
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -5039.5

$ws.Range("H27").Value = 4449.5
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4449.5
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4449.5
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -4663.5

$ws.Range("H55").Value = 1183.7273
$ws.Range("I55").Value = 910.44446
$ws.Range("J55").Value = 2413.5
$ws.Range("K55").Value = 910.44446
$ws.Range("L55").Value = 2413.5
$ws.Range("M55").Value = -737.44446
$ws.Range("N55").Value = -2759.5

$ws.Range("H93").Value = 1217.75
$ws.Range("I93").Value = 548.4
$ws.Range("J93").Value = 2333.3333
$ws.Range("K93").Value = 548.4
$ws.Range("L93").Value = 2333.3333
$ws.Range("M93").Value = 699.6
$ws.Range("N93").Value = -4829.3333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4250.5
$ws.Range("I81").Value = 4166.6665
$ws.Range("J81").Value = 4334.3335
$ws.Range("K81").Value = 8333.333000000001
$ws.Range("L81").Value = 8668.666999999999
$ws.Range("M81").Value = -7272.333000000001
$ws.Range("N81").Value = -10790.667

$ws.Range("H84").Value = 4250.5
$ws.Range("I84").Value = 4166.6665
$ws.Range("J84").Value = 4334.3335
$ws.Range("K84").Value = 41666.665
$ws.Range("L84").Value = 43343.335
$ws.Range("M84").Value = -36362.665
$ws.Range("N84").Value = -53951.335

$ws.Range("H105").Value = 15000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 15000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -21988

$ws.Range("H126").Value = 206998
$ws.Range("I126").Value = 502497.5
$ws.Range("J126").Value = 9998.333000000001
$ws.Range("K126").Value = 1507492.5
$ws.Range("L126").Value = 29994.999
$ws.Range("M126").Value = -1505022.5
$ws.Range("N126").Value = -34934.999
